$d = $word.ActiveDocument

# 1. Delete the "en instans ba af BeregnAfskrivning eksisterer" paragraph (Preconditions section)
$d.Paragraphs(14).Range.Delete()

# 2. Replace the whole Postconditions block (paragraphs 15-23) with the new
#    set of paragraphs describing the if/else flow.
$startP = $d.Paragraphs(15)
$endP = $d.Paragraphs(23)
$rng = $d.Range($startP.Range.Start, $endP.Range.End)

$TAB = [char]9
$CR = [char]13

$newText = (
    "Hvis en instans a af Afskrivning findes, hvor a.navn = navn og a er kendt af h" + $CR +
    $TAB + "a indeholder en instans sa af Saldoafskrivning" + $CR +
    $TAB + "sa.afskrivningsprocent blev sat til afskrivningsprocent" + $CR +
    $TAB + "a indeholder en instans av af Anskaffelsesværdi" + $CR +
    $TAB + "av.beløb blev sat til anskaffelsesværdi" + $CR +
    $TAB + "a.afskrivningsværdi blev sat til av.beløb * sa.afskrivningsprocent" + $CR +
    $TAB + "a blev præsenteret for h" + $CR +
    "Ellers" + $CR +
    $TAB + "En instans sa af Saldoafskrivning blev skabt" + $CR +
    $TAB + "sa.afskrivningsprocent blev sat til afskrivningsprocent" + $CR +
    $TAB + "En instans av af Anskaffelsesværdi blev skabt" + $CR +
    $TAB + "av.beløb blev sat til anskaffelsesværdi" + $CR +
    $TAB + "En instans a af Afskrivning blev skabt" + $CR +
    $TAB + "a blev sat til at indeholde sa" + $CR +
    $TAB + "a blev sat til at indeholde av" + $CR +
    $TAB + "a.navn blev sat til navn" + $CR +
    $TAB + "a.afskrivningsværdi blev sat til av.beløb * sa.afskrivningsprocent" + $CR +
    $TAB + "h blev sat til at kende a" + $CR +
    $TAB + "a blev præsenteret for h"
)

$rng.Text = $newText

# 3. Remove the trailing empty paragraph at the end of the document body by
#    merging it into the previous (now final) paragraph.
$count = $d.Paragraphs.Count
$secondLastEnd = $d.Paragraphs($count - 1).Range.End
$lastEnd = $d.Paragraphs($count).Range.End
$d.Range($secondLastEnd - 1, $lastEnd).Delete()

# 4. Re-add the _GoBack bookmark at the end of the final paragraph, matching
#    its original placement in the source document.
$finalP = $d.Paragraphs($d.Paragraphs.Count)
$d.Bookmarks.Add("_GoBack", $finalP.Range)

$i=1
foreach ($pp in $d.Paragraphs) {
  Write-Host $i "|" $pp.Range.Text
  $i=$i+1
}
